$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "#04-파이썬(Python) 문자열(str)"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-04"

$ws.Range("D9").Value = "지원자 Essay 관련 정보"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/math-stat-engineer-applicants/#utm_source=rss&utm_medium=rss&utm_campaign=math-stat-engineer-applicants"

$ws.Range("D27").Value = "TensorFlow Custom Op으로 데이터 변환 최적화하기"
$ws.Range("E27").Value = "https://blog.pingpong.us/custom-op-in-data-processing/"

$ws.Range("D28").Value = "Complementary Filter(LPF+HPF) Gyroscope, Accelerometer fusion"
$ws.Range("E28").Value = "https://ropiens.tistory.com/127"

$ws.Range("D29").Value = "프로메디우스"

$ws.Range("D32").Value = "RECSIM: A Conﬁgurable Simulation Platform for Recommender System 기초 개념"
$ws.Range("E32").Value = "https://dodonam.tistory.com/319"

$ws.Range("D51").Value = "[python] datetime 모듈로 일, 시간, 분, 초 더하거나 빼는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/1191"
